$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 5953939.5
$ws.Range("I19").Value = 13889528
$ws.Range("J19").Value = 2248.0833
$ws.Range("K19").Value = 13889528
$ws.Range("L19").Value = 2248.0833
$ws.Range("M19").Value = -13889353
$ws.Range("N19").Value = -2598.0833
$ws.Range("H32").Value = 1586
$ws.Range("I32").Value = 2900
$ws.Range("J32").Value = 1060.4
$ws.Range("K32").Value = 2900
$ws.Range("L32").Value = 1060.4
$ws.Range("M32").Value = -2574
$ws.Range("N32").Value = -1712.4
$ws.Range("H33").Value = 991.64703
$ws.Range("I33").Value = 1469
$ws.Range("J33").Value = 116.5
$ws.Range("K33").Value = 1469
$ws.Range("L33").Value = 116.5
$ws.Range("M33").Value = -1240
$ws.Range("N33").Value = -574.5
$ws.Range("H98").Value = 20927352
$ws.Range("I98").Value = 7144037
$ws.Range("J98").Value = 85249490
$ws.Range("K98").Value = 7144037
$ws.Range("L98").Value = 85249490
$ws.Range("M98").Value = -7142539
$ws.Range("N98").Value = -85252486
$ws.Range("H116").Value = 5317635
$ws.Range("I116").Value = 2860155
$ws.Range("J116").Value = 8389485
$ws.Range("K116").Value = 2860155
$ws.Range("L116").Value = 8389485
$ws.Range("M116").Value = -2856713
$ws.Range("N116").Value = -8396369
$ws.Range("H122").Value = 20927352
$ws.Range("I122").Value = 7144037
$ws.Range("J122").Value = 85249490
$ws.Range("K122").Value = 21432111
$ws.Range("L122").Value = 255748470
$ws.Range("M122").Value = -21429661
$ws.Range("N122").Value = -255753370
$ws.Range("H132").Value = 1974542.9
$ws.Range("I132").Value = 418649.06
$ws.Range("J132").Value = 12347168
$ws.Range("K132").Value = 1255947.18
$ws.Range("L132").Value = 37041504
$ws.Range("M132").Value = -1253417.18
$ws.Range("N132").Value = -37046564
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 6588
$ws.Range("I2").Value = 7992.638
$ws.Range("J2").Value = 586.36365
$ws.Range("K2").Value = 7992.638
$ws.Range("L2").Value = 586.36365
$ws.Range("M2").Value = -7879.638
$ws.Range("N2").Value = -812.36365
$ws.Range("H61").Value = 4188139.5
$ws.Range("I61").Value = 1894855.5
$ws.Range("J61").Value = 29414264
$ws.Range("K61").Value = 1894855.5
$ws.Range("L61").Value = 29414264
$ws.Range("M61").Value = -1894643.5
$ws.Range("N61").Value = -29414688
$ws.Range("H63").Value = 1590.2703
$ws.Range("I63").Value = 1545.2941
$ws.Range("J63").Value = 2100
$ws.Range("K63").Value = 1545.2941
$ws.Range("L63").Value = 2100
$ws.Range("M63").Value = -859.2941000000001
$ws.Range("N63").Value = -3472
$ws.Range("H66").Value = 1590.2703
$ws.Range("I66").Value = 1545.2941
$ws.Range("J66").Value = 2100
$ws.Range("K66").Value = 7726.4705
$ws.Range("L66").Value = 10500
$ws.Range("M66").Value = -4294.4705
$ws.Range("N66").Value = -17364
$ws.Range("H74").Value = 45099664
$ws.Range("I74").Value = 40000692
$ws.Range("J74").Value = 59263470
$ws.Range("K74").Value = 40000692
$ws.Range("L74").Value = 59263470
$ws.Range("M74").Value = -39999818
$ws.Range("N74").Value = -59265218
$ws.Range("H77").Value = 45099664
$ws.Range("I77").Value = 40000692
$ws.Range("J77").Value = 59263470
$ws.Range("K77").Value = 200003460
$ws.Range("L77").Value = 296317350
$ws.Range("M77").Value = -199999092
$ws.Range("N77").Value = -296326086
$ws.Range("H97").Value = 490.375
$ws.Range("I97").Value = 475.15
$ws.Range("J97").Value = 515.75
$ws.Range("K97").Value = 475.15
$ws.Range("L97").Value = 515.75
$ws.Range("M97").Value = 20.85000000000002
$ws.Range("N97").Value = -1507.75
$ws.Range("H116").Value = 6588
$ws.Range("I116").Value = 7992.638
$ws.Range("J116").Value = 586.36365
$ws.Range("K116").Value = 7992.638
$ws.Range("L116").Value = 586.36365
$ws.Range("M116").Value = -5698.638
$ws.Range("N116").Value = -5174.36365
$ws.Range("H122").Value = 1194.6451
$ws.Range("I122").Value = 1196.6666
$ws.Range("J122").Value = 1187.7142
$ws.Range("K122").Value = 3589.9998
$ws.Range("L122").Value = 3563.1426
$ws.Range("M122").Value = -1139.9998
$ws.Range("N122").Value = -8463.142599999999
$ws.Range("H132").Value = 21988762
$ws.Range("I132").Value = 24330282
$ws.Range("J132").Value = 7939650.5
$ws.Range("K132").Value = 72990846
$ws.Range("L132").Value = 23818951.5
$ws.Range("M132").Value = -72988316
$ws.Range("N132").Value = -23824011.5
$ws.Range("H136").Value = 4188139.5
$ws.Range("I136").Value = 1894855.5
$ws.Range("J136").Value = 29414264
$ws.Range("K136").Value = 5684566.5
$ws.Range("L136").Value = 88242792
$ws.Range("M136").Value = -5682016.5
$ws.Range("N136").Value = -88247892
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 6588
$ws.Range("I3").Value = 7992.638
$ws.Range("J3").Value = 586.36365
$ws.Range("K3").Value = 7992.638
$ws.Range("L3").Value = 586.36365
$ws.Range("M3").Value = -7878.638
$ws.Range("N3").Value = -814.36365
$ws.Range("H105").Value = 1792.8572
$ws.Range("I105").Value = 1861.5385
$ws.Range("J105").Value = 900
$ws.Range("K105").Value = 1861.5385
$ws.Range("L105").Value = 900
$ws.Range("M105").Value = -114.5385000000001
$ws.Range("N105").Value = -4394
$ws.Range("H107").Value = 794.65
$ws.Range("I107").Value = 735
$ws.Range("J107").Value = 1132.6666
$ws.Range("K107").Value = 735
$ws.Range("L107").Value = 1132.6666
$ws.Range("M107").Value = 1185
$ws.Range("N107").Value = -4972.6666
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 695.75
$ws.Range("I22").Value = 204.33333
$ws.Range("J22").Value = 2170
$ws.Range("K22").Value = 204.33333
$ws.Range("L22").Value = 2170
$ws.Range("M22").Value = 145.66667
$ws.Range("N22").Value = -2870
$ws.Range("H31").Value = 1781.62
$ws.Range("I31").Value = 661.807
$ws.Range("J31").Value = 3266.0232
$ws.Range("K31").Value = 661.807
$ws.Range("L31").Value = 3266.0232
$ws.Range("M31").Value = -366.807
$ws.Range("N31").Value = -3856.0232
$ws.Range("H34").Value = 1781.62
$ws.Range("I34").Value = 661.807
$ws.Range("J34").Value = 3266.0232
$ws.Range("K34").Value = 661.807
$ws.Range("L34").Value = 3266.0232
$ws.Range("M34").Value = -459.807
$ws.Range("N34").Value = -3670.0232
$ws.Range("H105").Value = 7457.3335
$ws.Range("I105").Value = 1736.125
$ws.Range("J105").Value = 18899.75
$ws.Range("K105").Value = 1736.125
$ws.Range("L105").Value = 18899.75
$ws.Range("M105").Value = 10.875
$ws.Range("N105").Value = -22393.75
$ws.Range("H107").Value = 772.27026
$ws.Range("I107").Value = 268.5
$ws.Range("J107").Value = 869.7742
$ws.Range("K107").Value = 268.5
$ws.Range("L107").Value = 869.7742
$ws.Range("M107").Value = 1651.5
$ws.Range("N107").Value = -4709.7742
$ws.Range("H122").Value = 3713.3555
$ws.Range("I122").Value = 4894.8623
$ws.Range("K122").Value = 14684.5869
$ws.Range("M122").Value = -12234.5869
$ws.Range("H132").Value = 1996.8918
$ws.Range("J132").Value = 3723.5557
$ws.Range("L132").Value = 11170.6671
$ws.Range("N132").Value = -16230.6671
$ws.Range("H134").Value = 1434446.2
$ws.Range("I134").Value = 5259.2
$ws.Range("J134").Value = 13344338
$ws.Range("K134").Value = 15777.6
$ws.Range("L134").Value = 40033014
$ws.Range("M134").Value = -13242.6
$ws.Range("N134").Value = -40038084
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 4879649
$ws.Range("J92").Value = 6099311
$ws.Range("L92").Value = 18297933
$ws.Range("N92").Value = -18300429
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 5806.6665
$ws.Range("I102").Value = 5901.6665
$ws.Range("J102").Value = 5331.6665
$ws.Range("K102").Value = 5901.6665
$ws.Range("L102").Value = 5331.6665
$ws.Range("M102").Value = -4279.6665
$ws.Range("N102").Value = -8575.666499999999
$ws.Range("H107").Value = 195.33333
$ws.Range("I107").Value = 194.61539
$ws.Range("J107").Value = 200
$ws.Range("K107").Value = 194.61539
$ws.Range("L107").Value = 200
$ws.Range("M107").Value = 1725.38461
$ws.Range("N107").Value = -4040
$ws.Range("H122").Value = 5210501.5
$ws.Range("I122").Value = 2006.88
$ws.Range("J122").Value = 23812268
$ws.Range("K122").Value = 6020.64
$ws.Range("L122").Value = 71436804
$ws.Range("M122").Value = -3570.64
$ws.Range("N122").Value = -71441704
$ws.Range("H126").Value = 6110.919
$ws.Range("I126").Value = 14275.5
$ws.Range("J126").Value = 3858.6206
$ws.Range("K126").Value = 42826.5
$ws.Range("L126").Value = 11575.8618
$ws.Range("M126").Value = -40356.5
$ws.Range("N126").Value = -16515.8618
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 100000920
$ws.Range("I46").Value = 1249.5
$ws.Range("J46").Value = 166667360
$ws.Range("K46").Value = 1249.5
$ws.Range("L46").Value = 166667360
$ws.Range("M46").Value = -1061.5
$ws.Range("N46").Value = -166667736
$ws.Range("H132").Value = 5960146.5
$ws.Range("I132").Value = 7525343
$ws.Range("J132").Value = 12399.8
$ws.Range("K132").Value = 22576029
$ws.Range("L132").Value = 37199.39999999999
$ws.Range("M132").Value = -22573499
$ws.Range("N132").Value = -42259.39999999999
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 6810.125
$ws.Range("I100").Value = 8916.833000000001
$ws.Range("J100").Value = 490
$ws.Range("K100").Value = 17833.666
$ws.Range("L100").Value = 980
$ws.Range("M100").Value = -17292.666
$ws.Range("N100").Value = -2062
$ws.Range("H132").Value = 692289.9
$ws.Range("I132").Value = 1770.4773
$ws.Range("K132").Value = 5311.4319
$ws.Range("M132").Value = -2781.4319
$ws.Range("H136").Value = 1202.1136
$ws.Range("I136").Value = 647.18866
$ws.Range("J136").Value = 2042.4286
$ws.Range("K136").Value = 1941.56598
$ws.Range("L136").Value = 6127.2858
$ws.Range("M136").Value = 608.4340199999999
$ws.Range("N136").Value = -11227.2858
